# Adição da Sprint 9 e Alterações no Burndown
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Update the "HORAS TRABALHADAS" (burndown) column values for the
#    existing rows (column F).
# ---------------------------------------------------------------------
$horas = @{
    3  = " 3 Horas e 15 Minutos"
    4  = " 3 Horas e 15 Minutos"
    6  = "2 Horas 30 Minutos"
    7  = "5 Horas"
    8  = "1 Hora e 30 Minutos"
    10 = "2 Horas"
    11 = "2 Horas"
    12 = "2 Horas"
    14 = "2 Horas"
    15 = "2 Horas e 20 Minutos"
    16 = "2 Hora e 10 Minutos"
    18 = "2 Horas"
    19 = "3 Hora e 30 Minutos"
    20 = "2 Hora e 30 Minutos"
    22 = "2 Horas"
    23 = "  4 Horas e 30 Minutos"
    24 = " 1 Hora e 30 Minutos"
    25 = " 1 Hora"
    27 = "1 Hora 30 Minutos"
    28 = "2 Horas"
    29 = "2 Horas e 30 Minutos"
    30 = "1 Horas e 30 Minutos"
    31 = "1 Horas e 30 Minutos"
    33 = "1 Hora"
    34 = "1 Hora"
    35 = "1 Hora"
    36 = "1 Hora e 30 Minutos"
    38 = "1 Hora"
    39 = "1 Hora"
}

foreach ($row in $horas.Keys) {
    $ws.Range("F$row").Value = $horas[$row]
}

# F8 additionally picks up a distinct (font-reasserted) style in the
# source workbook - reassert the theme color so a new cellXfs entry is
# produced, mirroring the authored change.
$ws.Range("F8").Font.ThemeColor = 1

# ---------------------------------------------------------------------
# 2) Resize the table / used range to make room for the new Sprint 9
#    rows (was B2:G39, becomes B2:G44) and add a separator row (40)
#    plus four Sprint 9 task rows (41-44).
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:G44"))

# Copy formatting from the existing blank separator row down to row 40.
$ws.Range("B32:G32").Copy()
$ws.Range("B40:G40").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Copy formatting from a clean data block (rows 33-36) down to the new
# Sprint 9 rows (41-44).
$ws.Range("B33:G36").Copy()
$ws.Range("B41:G44").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ---------------------------------------------------------------------
# 3) Fill in the Sprint 9 data.
# ---------------------------------------------------------------------
$ws.Range("B41:B44").Value = "Sprint 9"
$ws.Range("C41:C44").Value = "15/05/2024 - 22/05/2024"
$ws.Range("F41:F44").Value = "30 Minutos"
$ws.Range("G41:G44").Value = "Pronto"

$ws.Range("D41").Value = "Mudanças finais no protótipo"
$ws.Range("E41").Value = "Cristielen"

$ws.Range("D42").Value = "Edições nos Diagramas"
$ws.Range("E42").Value = "Guilherme"

$ws.Range("D43").Value = "Edições na planilha de testes"
$ws.Range("E43").Value = "Nicolas"

$ws.Range("D44").Value = "Finalização do Manual do Usuário"
$ws.Range("E44").Value = "Bruno"

# ---------------------------------------------------------------------
# 4) View tweaks captured in the diff: gridlines hidden and a new
#    active selection.
# ---------------------------------------------------------------------
$aw = $excel.ActiveWindow()
$aw.DisplayGridlines = $false
$ws.Range("L28").Select()

Write-Host "Sprint 9 added and burndown column refreshed."
